# HermitsRestDuplicateSlugImport.xlsx fix:
# The "Longitude" calculated column in Table13 was wired to the wrong
# source column (Latitude, column J) instead of the actual longitude
# values in column K. Re-point the column's formula at K2 (the table
# will propagate the relative formula down every data row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$table = $ws.ListObjects.Item("Table13")
$longitudeColumn = $table.ListColumns.Item("Longitude")
$longitudeColumn.DataBodyRange.Formula = "=K2"

# Leave the same cell selected that was active when the fix was made.
$ws.Range("AC2").Select()
